# "add perimeter to UI"
# Appends a new block of benchmark rows (perimeter/recall-style summary block)
# to the bottom of the results sheet, mirroring the layout already used
# throughout the sheet for the other model/K-field summary blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-PlainCell($row, $col, $text) {
    # Leaves the cell on the column's inherited/default style (matches the
    # vast majority of cells in the new rows, which just pick up style 10
    # for columns B:G, or style 1 for column A).
    $ws.Cells.Item($row, $col).Value = $text
}

function Set-NormalCell($row, $col, $text) {
    # Forces the "Normal" (un-styled / non-bold Arial 10) look for cells
    # that must NOT inherit the column's default style.
    $c = $ws.Cells.Item($row, $col)
    $c.Value = $text
    $c.Style = "Normal"
}

# ---------------------------------------------------------------------
# Row 81 / 82 : att_unet ... Kfield/3 and Kfield/4 summary rows
# ---------------------------------------------------------------------
Set-NormalCell 81 1 "att_unet-diana_healthy_marfan-lr_0.001-batch_2-augmented-instance_normalization-polygon2mask-Kfield/3"
Set-PlainCell  81 3 "74.13% (+/- 15.17%)"
Set-PlainCell  81 4 "60.94% (+/- 17.19%)"
Set-PlainCell  81 5 "16.78 (+/- 7.57)"
Set-PlainCell  81 6 "80.38% (+/- 16.93%)"
Set-PlainCell  81 7 " 81.63% (+/- 17.35%)"
$ws.Rows(81).RowHeight = 12.5

Set-NormalCell 82 1 "att_unet-diana_healthy_marfan-lr_0.001-batch_2-augmented-instance_normalization-polygon2mask-Kfield/4"
Set-PlainCell  82 2 "loss  0.02% (+/- 0.02%)"
Set-PlainCell  82 3 "97.64% (+/- 1.61%)"
Set-PlainCell  82 4 "95.43% (+/- 2.85%)"
Set-PlainCell  82 5 "5.66 (+/- 2.01)"
Set-PlainCell  82 6 "98.62% (+/- 1.15%)"
Set-PlainCell  82 7 "96.83% (+/- 3.01%)"
$ws.Rows(82).RowHeight = 12.5

# ---------------------------------------------------------------------
# Row 87 : new dataset header / metric row ("6480 - noisy_waves")
# ---------------------------------------------------------------------
Set-PlainCell  87 1 "6480 - noisy_waves"
Set-NormalCell 87 2 "loss  0.04% (+/- 0.14%)"
Set-PlainCell  87 3 "dice_coef  95.88% (+/- 14.00%)"
Set-PlainCell  87 4 "iou  94.03% (+/- 14.54%)"
Set-PlainCell  87 5 "hausdorff  6.16 (+/- 6.05)"
Set-PlainCell  87 6 "precision  95.95% (+/- 14.65%)"
Set-PlainCell  87 7 "recall  96.05% (+/- 12.61%)"

# ---------------------------------------------------------------------
# Row 90 / 91 : att_res_unet ... Kfield/3 and Kfield/4 summary rows
# ---------------------------------------------------------------------
Set-NormalCell 90 1 "att_res_unet-diana_healthy_marfan-lr_0.001-batch_2-augmented-instance_normalization-polygon2mask-Kfield/3"
Set-PlainCell  90 2 "0.02% (+/- 0.01%)"
Set-PlainCell  90 3 "98.03% (+/- 0.78%)"
Set-PlainCell  90 4 "96.14% (+/- 1.47%)"
Set-PlainCell  90 5 "5.24 (+/- 1.72)"
Set-PlainCell  90 6 "98.14% (+/- 1.27%)"
Set-PlainCell  90 7 "97.99% (+/- 1.64%)"
$ws.Rows(90).RowHeight = 12.5

Set-NormalCell 91 1 "att_res_unet-diana_healthy_marfan-lr_0.001-batch_2-augmented-instance_normalization-polygon2mask-Kfield/4"
Set-PlainCell  91 2 "loss  0.02% (+/- 0.01%)"
Set-PlainCell  91 3 "dice_coef  97.95% (+/- 0.76%)"
Set-PlainCell  91 4 "iou  96.00% (+/- 1.44%)"
Set-PlainCell  91 5 "hausdorff  5.28 (+/- 1.70)"
Set-PlainCell  91 6 "precision_1  98.12% (+/- 1.47%)"
Set-PlainCell  91 7 "recall_1  97.87% (+/- 1.39%)"
$ws.Rows(91).RowHeight = 12.5

# ---------------------------------------------------------------------
# Column A grew wider to fit the new (longer) run names.
# ---------------------------------------------------------------------
$ws.Columns("A:A").ColumnWidth = 91.90625

# ---------------------------------------------------------------------
# Leave the view pointed at the newly-added block, like the author did.
# ---------------------------------------------------------------------
$ws.Range("C90:G90").Select() | Out-Null
